$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "12:14:59"
$ws.Range("C2").Value = "12:15:00"
$ws.Range("D2").Value = "12:15:01"
$ws.Range("E2").Value = "12:15:02"
